$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2026-01-22"
$ws.Cells.Item($row, 2).Value = "2026-01-22 21:30:41"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Value = 1.32
$ws.Cells.Item($row, 4).Value = 4
$ws.Cells.Item($row, 5).Value = 0.57
$ws.Cells.Item($row, 6).Value = 21
$ws.Cells.Item($row, 7).Value = 0.7500000000000001
$ws.Cells.Item($row, 8).Value = 25
$ws.Cells.Item($row, 9).Value = 2.2
$ws.Cells.Item($row, 10).Value = 129
$ws.Cells.Item($row, 11).Value = 1.1
$ws.Cells.Item($row, 12).Value = 1658
$ws.Cells.Item($row, 13).Value = 1.1
$ws.Cells.Item($row, 14).Value = 1787
